$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.501.26"
$ws.Range("E2").Value = "  -0.02%  "

# Row 3
$ws.Range("D3").Value = "1.734.14"
$ws.Range("E3").Value = "  -0.13%  "

# Row 4
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").Value = "'246.86"
$ws.Range("E5").Value = "  +1.30%  "

# Row 6
$ws.Range("E6").Value = "  -0.05%  "

# Row 7
$ws.Range("D7").Value = "'0.4887"
$ws.Range("E7").Value = "  +1.94%  "

# Row 8
$ws.Range("D8").Value = "'0.2666"
$ws.Range("E8").Value = "  +0.10%  "

# Row 9
$ws.Range("D9").Value = "'0.06306"
$ws.Range("E9").Value = "  +1.32%  "

# Row 10
$ws.Range("D10").Value = "1.729.29"
$ws.Range("E10").Value = "  -0.46%  "

# Row 11
$ws.Range("D11").Value = "'0.07025"
$ws.Range("E11").Value = "  -1.50%  "

# Row 12
$ws.Range("D12").Value = "'15.68"
$ws.Range("E12").Value = "  -0.29%  "

# Row 13
$ws.Range("D13").Value = "'4.594"
$ws.Range("E13").Value = "  +1.39%  "

# Row 14
$ws.Range("D14").Value = "'0.6083"
$ws.Range("E14").Value = "  -1.24%  "

# Row 15
$ws.Range("D15").Value = "'77.34"
$ws.Range("E15").Value = "  +0.67%  "

# Row 16
$ws.Range("E16").Value = "  -0.07%  "

# Row 17
$ws.Range("D17").Value = "'0.000007504"
$ws.Range("E17").Value = "  +8.95%  "

# Row 18
$ws.Range("D18").Value = "26.489.37"
$ws.Range("E18").Value = "  -0.11%  "

# Row 19
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  -0.16%  "

# Row 20
$ws.Range("D20").Value = "'11.52"
$ws.Range("E20").Value = "  -1.78%  "

# Row 21
$ws.Range("D21").Value = "1.951.86"
$ws.Range("E21").Value = "  -0.46%  "

# Row 22
$ws.Range("D22").Value = "'4.576"
$ws.Range("E22").Value = "  +0.36%  "

# Row 23
$ws.Range("D23").Value = "'8.704"
$ws.Range("E23").Value = "  -2.05%  "

# Row 24
$ws.Range("D24").Value = "'5.226"
$ws.Range("E24").Value = "  -1.99%  "

# Row 25
$ws.Range("D25").Value = "'140.69"
$ws.Range("E25").Value = "  +3.77%  "

# Row 26
$ws.Range("E26").Value = "  +0.55%  "

# Row 27
$ws.Range("D27").Value = "'1.416"
$ws.Range("E27").Value = "  +0.48%  "

# Row 28
$ws.Range("D28").Value = "'1.767"
$ws.Range("E28").Value = "  -1.77%  "

# Row 29
$ws.Range("D29").Value = "'107.95"
$ws.Range("E29").Value = "  +1.32%  "

# Row 30
$ws.Range("D30").Value = "'4.033"
$ws.Range("E30").Value = "  +1.14%  "

# Row 31
$ws.Range("D31").Value = "'0.08006"
$ws.Range("E31").Value = "  +1.56%  "

# Row 32
$ws.Range("D32").Value = "'3.712"
$ws.Range("E32").Value = "  +0.05%  "

# Row 33
$ws.Range("D33").Value = "'0.04574"
$ws.Range("E33").Value = "  +0.07%  "

# Row 34
$ws.Range("D34").Value = "'0.9999"
$ws.Range("E34").Value = "  -0.09%  "

# Row 35
$ws.Range("E35").Value = "  -0.19%  "

# Row 36
$ws.Range("D36").Value = "'1.008"
$ws.Range("E36").Value = "  +1.35%  "

# Row 37
$ws.Range("D37").Value = "'0.6351"
$ws.Range("E37").Value = "  +0.27%  "

# Row 38
$ws.Range("D38").Value = "'0.8937"
$ws.Range("E38").Value = "  -3.95%  "

# Row 39
$ws.Range("D39").Value = "'2.012"
$ws.Range("E39").Value = "  +1.84%  "

# Row 40
$ws.Range("D40").Value = "'2.397"
$ws.Range("E40").Value = "  -1.89%  "

# Row 41
$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = "  -0.30%  "

# Row 42
$ws.Range("D42").Value = "'0.01503"
$ws.Range("E42").Value = "  -0.29%  "

# Row 43
$ws.Range("D43").Value = "'101.71"
$ws.Range("E43").Value = "  -8.51%  "

# Row 44
$ws.Range("D44").Value = "'5.404"
$ws.Range("E44").Value = "  -5.08%  "

# Row 45
$ws.Range("D45").Value = "'0.3882"
$ws.Range("E45").Value = "  -0.35%  "

# Row 46
$ws.Range("D46").Value = "'6.913"
$ws.Range("E46").Value = "  +0.38%  "

# Row 47
$ws.Range("D47").Value = "'0.1183"
$ws.Range("E47").Value = "  -0.80%  "

# Row 48
$ws.Range("D48").Value = "'0.05391"
$ws.Range("E48").Value = "  +1.07%  "

# Row 49
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'30.46"
$ws.Range("E49").Value = "  -0.99%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.748"
$ws.Range("E50").Value = "  -1.51%  "

# Row 51
$ws.Range("E51").Value = "  +0.52%  "

